$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Through 2022-07-28"

# Update header text for July row
$ws.Range("A8").Value = "July (through 07-28)"

# Update July row values (row 8)
$ws.Range("C8").Value = 48
$ws.Range("D8").Value = 69
$ws.Range("E8").Value = 65
$ws.Range("F8").Value = 47
$ws.Range("G8").Value = 129
$ws.Range("H8").Value = 138
$ws.Range("I8").Value = 158

# Update Total row values (row 9)
$ws.Range("C9").Value = 296
$ws.Range("D9").Value = 459
$ws.Range("E9").Value = 418
$ws.Range("F9").Value = 298
$ws.Range("G9").Value = 601
$ws.Range("H9").Value = 898
$ws.Range("I9").Value = 964
